# Updates "苏州-漫展信息" workbook to the newly scraped bilibili show data.
# - Refreshes a handful of "interested count" (F column) numbers on existing
#   rows (both on "展览" and the merged "全部类型" sheet).
# - Inserts one newly-discovered event ("苏州·首届Redamancy动漫游戏嘉年华",
#   2024.04.20) ahead of the existing chronological rows, pushing every
#   following row down by one and renumbering the leading index column.

function Set-EventRow {
    param(
        $ws, $r,
        $idx, $start, $name, $place,
        $range, $want, $price, $link, $cover
    )
    $ws.Cells.Item($r, 1).Value = $idx
    # Column B holds plain "YYYY.MM.DD" text in the source data; without
    # forcing a text format first, Excel auto-coerces it into a date serial.
    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $start
    $ws.Cells.Item($r, 3).Value = $name
    $ws.Cells.Item($r, 4).Value = $place
    $ws.Cells.Item($r, 5).Value = $range
    $ws.Cells.Item($r, 6).Value = $want
    $ws.Cells.Item($r, 7).Value = $price
    $ws.Cells.Item($r, 8).Value = $link
    $ws.Cells.Item($r, 9).Value = $cover
}

function Update-ExhibitionSheet {
    param($ws, $firstNewRow)

    # Insert a fresh row above the block that is about to be renumbered, and
    # restore the bordered/centered index-column formatting that "Insert"
    # does not fully carry over from the row above.
    $ws.Rows.Item($firstNewRow).Insert()
    $copySrc = "A" + ($firstNewRow + 1)
    $copyDst = "A" + $firstNewRow
    $ws.Range($copySrc).Copy()
    $ws.Range($copyDst).PasteSpecial(-4122)
    $excel.CutCopyMode = 0

    $r = $firstNewRow
    $idx = $r - 1
    Set-EventRow $ws $r $idx "2024.04.20" "苏州·首届Redamancy动漫游戏嘉年华" "清禾路886号 尹山湖大剧院" `
        "2024.04.20 10:00-04.20 17:00" 0 60 `
        "https://show.bilibili.com/platform/detail.html?id=81879" `
        "//i0.hdslb.com/bfs/openplatform/202402/lR4oJWzI1708309129629.jpeg"

    $r = $firstNewRow + 1
    $idx = $r - 1
    Set-EventRow $ws $r $idx "2024.04.21" "苏州·梦幻岛 国乙主题文化展（日夜场） 梦幻岛之约3.0" "清禾路888号2号楼3楼 格莱美婚礼宴会中心" `
        "2024.04.21 10:00-04.21 21:00" 593 49.9 `
        "https://show.bilibili.com/platform/detail.html?id=78666" `
        "//i0.hdslb.com/bfs/openplatform/202312/X0PZ3YhH1703822037665.jpeg"

    $r = $firstNewRow + 2
    $idx = $r - 1
    Set-EventRow $ws $r $idx "2024.05.01" "昆山·第十二届理想乡动漫游戏展" "花桥经济开发区绿地大道1598号 花桥国际博览中心" `
        "2024.05.01 10:00-05.03 17:00" 11674 59 `
        "https://show.bilibili.com/platform/detail.html?id=77196" `
        "//i2.hdslb.com/bfs/openplatform/202310/9xMTQMlg1696736126094.png"

    $r = $firstNewRow + 3
    $idx = $r - 1
    Set-EventRow $ws $r $idx "2024.05.01" "苏州·第十七届 I COME ACG  动漫品牌博览会" "金山南路288号 广电国际会展中心" `
        "2024.05.01 10:00-05.02 17:00" 11962 65 `
        "https://show.bilibili.com/platform/detail.html?id=79789" `
        "//i2.hdslb.com/bfs/openplatform/202312/lau3mW031702535438289.jpeg"

    $r = $firstNewRow + 4
    $idx = $r - 1
    Set-EventRow $ws $r $idx "2024.05.02" "昆山·第十二届理想乡动漫游戏展嘉宾北齐后主签售会" "花桥经济开发区绿地大道1598号 花桥国际博览中心" `
        "2024.05.02 14:00-05.02 16:00" 24 1 `
        "https://show.bilibili.com/platform/detail.html?id=81116" `
        "//i2.hdslb.com/bfs/openplatform/202401/EubrAneC1705648695005.jpeg"

    $r = $firstNewRow + 5
    $idx = $r - 1
    Set-EventRow $ws $r $idx "2024.05.02" "昆山·第十二届理想乡动漫游戏展嘉宾啊川签售会" "花桥经济开发区绿地大道1598号 花桥国际博览中心" `
        "2024.05.02 14:00-05.02 16:00" 84 1 `
        "https://show.bilibili.com/platform/detail.html?id=81100" `
        "//i2.hdslb.com/bfs/openplatform/202401/F24i5GMX1705646667852.jpeg"

    $r = $firstNewRow + 6
    $idx = $r - 1
    Set-EventRow $ws $r $idx "2024.05.02" "昆山·第十二届理想乡动漫游戏展嘉宾漠小然签售会" "花桥经济开发区绿地大道1598号 花桥国际博览中心" `
        "2024.05.02 14:00-05.02 16:00" 16 1 `
        "https://show.bilibili.com/platform/detail.html?id=81119" `
        "//i2.hdslb.com/bfs/openplatform/202401/SDnLB1gR1705648838683.jpeg"

    $r = $firstNewRow + 7
    $idx = $r - 1
    Set-EventRow $ws $r $idx "2024.05.02" "昆山·第十二届理想乡动漫游戏展嘉宾葫芦岛老八签售会" "花桥经济开发区绿地大道1598号 花桥国际博览中心" `
        "2024.05.02 14:00-05.02 16:00" 21 1 `
        "https://show.bilibili.com/platform/detail.html?id=81118" `
        "//i2.hdslb.com/bfs/openplatform/202401/uHOCneLv1705648779163.jpeg"

    $r = $firstNewRow + 8
    $idx = $r - 1
    Set-EventRow $ws $r $idx "2024.05.02" "昆山·第十二届理想乡动漫游戏展嘉宾沈辞签售会" "花桥经济开发区绿地大道1598号 花桥国际博览中心" `
        "2024.05.02 14:00-05.02 16:00" 64 1 `
        "https://show.bilibili.com/platform/detail.html?id=81120" `
        "//i0.hdslb.com/bfs/openplatform/202401/4Pay1rR61705648901961.jpeg"

    $r = $firstNewRow + 9
    $idx = $r - 1
    Set-EventRow $ws $r $idx "2024.05.03" "昆山·第十二届理想乡动漫游戏展嘉宾矮乐多aliga签售会" "花桥经济开发区绿地大道1598号 花桥国际博览中心" `
        "2024.05.03 14:00-05.03 16:00" 31 1 `
        "https://show.bilibili.com/platform/detail.html?id=81114" `
        "//i1.hdslb.com/bfs/openplatform/202401/Peub7FOc1705648580577.jpeg"
}

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: 展览 (exhibitions) ----
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F7").Value = 2655
$ws1.Range("F8").Value = 1163
$ws1.Range("F11").Value = 9754
$ws1.Range("F13").Value = 242
Update-ExhibitionSheet $ws1 14

# ---- Sheet 4: 全部类型 (all categories merged) ----
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F7").Value = 2655
$ws4.Range("F9").Value = 1163
$ws4.Range("F12").Value = 9754
$ws4.Range("F14").Value = 242
Update-ExhibitionSheet $ws4 15
